$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: swap the two halves of the resource-usage snapshot -----------
# Old: B2:E2 = 9,11,13,22   F2:I2 = 12,13,17,28
# New: B2:E2 = 12,13,17,28  F2:I2 = 9,11,13,22
$row2 = @(12, 13, 17, 28, 9, 11, 13, 22)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, 2 + $i).Value = $row2[$i]
}

# --- New rows: s3s4 totals, final totals, and a trailing grand total -----
$ws.Cells.Item(11, 1).Value = "s3s4"
$row11 = @(0, 0, 0, 0)
for ($i = 0; $i -lt $row11.Length; $i++) {
    $ws.Cells.Item(11, 2 + $i).Value = $row11[$i]
}

$ws.Cells.Item(12, 1).Value = "final"
$row12 = @(2, 9, 10, 16)
for ($i = 0; $i -lt $row12.Length; $i++) {
    $ws.Cells.Item(12, 2 + $i).Value = $row12[$i]
}

$ws.Cells.Item(13, 5).Value = 107

# --- Sheet view: scroll position + selection ------------------------------
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("I19").Select() | Out-Null
